$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.643.63'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.643.41'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.61'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.17%  '
$ws.Range('E6').Value = '  +1.00%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('E9').Value = '  +0.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.25'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.07%  '
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.872.72'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.73%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.21'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.94%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.638.70'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.20%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.531'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.69'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.59%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.666.65'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0750'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '217.14'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.89%  '
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('E21').Value = '  +1.66%  '
$ws.Range('E22').Value = '  +2.15%  '
$ws.Range('E23').Value = '  +2.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +9.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.12'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('E27').Value = '  -1.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.13'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.05%  '
$ws.Range('E29').Value = '  +1.34%  '
$ws.Range('E30').Value = '  +2.56%  '
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('E32').Value = '  +2.61%  '
$ws.Range('E33').Value = '  +1.76%  '
$ws.Range('E34').Value = '  +2.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.270.10'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.43%  '
$ws.Range('E36').Value = '  +5.11%  '
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('E38').Value = '  +5.90%  '
$ws.Range('E39').Value = '  +2.67%  '
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.811'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.35%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.47'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.56%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.24'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.782.86'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '93.15'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.46'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.22%  '
$ws.Range('E47').Value = '  +3.09%  '
$ws.Range('E48').Value = '  +0.99%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.77'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.28%  '
$ws.Range('E50').Value = '  +3.45%  '
$ws.Range('E51').Value = '  -0.49%  '
